$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '29.922.85'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.876.54'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7429'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3150'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07209'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08386'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('D12').Value = '1.908.40'
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7523'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.415'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.64'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.19%  '
$ws.Range('D16').Value = '29.958.55'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.074'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.66%  '
$ws.Range('E19').Value = '  -1.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007859'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').Value = '2.129.67'
$ws.Range('E22').Value = '  -4.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.048'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -5.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.272'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.520'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.610'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.536'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.287'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05331'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.237'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7494'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.002'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01966'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.759'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4535'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('D42').Value = '1.112.85'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.086'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.41'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8561'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.618'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').Value = '2.028.67'
$ws.Range('E50').Value = '  -3.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.910'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.42%  '
